$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43/44 swap (HuobiToken <-> VeChain) ---
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Range("D43").Value = "0.0218"
$ws.Range("E43").Value = "  +1.36%  "

$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Range("D44").Value = "2.90"
$ws.Range("E44").Value = "  -1.24%  "

# --- Price (D) cells that need Text format to avoid numeric auto-conversion ---
$numFmtRows = @(5,6,7,9,10,11,13,15,17,20,22,23,25,26,27,28,29,32,34,35,36,37,40,45,46,47,48,50)
foreach ($r in $numFmtRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# --- Apply D and E values per row ---
$ws.Range("D2").Value = "36.618.06"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.012.26"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "247.25"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("D7").Value = "63.05"
$ws.Range("E7").Value = "  +2.59%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.386"
$ws.Range("E9").Value = "  +4.29%  "
$ws.Range("D10").Value = "57.19"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  +6.85%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "0.886"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("E14").Value = "  +13.50%  "
$ws.Range("D15").Value = "14.21"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").Value = "2.308.04"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "5.56"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "2.012.67"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "36.538.40"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "71.96"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("D22").Value = "5.36"
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").Value = "238.63"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  -8.85%  "
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").Value = "10.07"
$ws.Range("E27").Value = "  +5.47%  "
$ws.Range("D28").Value = "159.42"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").Value = "20.19"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("E30").Value = "  +18.53%  "
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").Value = "5.02"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "0.0634"
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("D35").Value = "4.52"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "6.46"
$ws.Range("E36").Value = "  +9.47%  "
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "3.19"
$ws.Range("E40").Value = "  +15.00%  "
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D45").Value = "1.13"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "96.16"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").Value = "16.76"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "7.70"
$ws.Range("E48").Value = "  -4.34%  "
$ws.Range("D49").Value = "1.364.04"
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "2.199.62"
$ws.Range("E51").Value = "  +0.40%  "
